$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 98 (pushes former rows 98-110 down to 99-111)
$ws.Rows.Item(98).Insert()

# Populate the new row 98 with a new weekly price record
# (same market/product metadata as its neighbours, new date + origin)
$ws.Cells.Item(98, 1).Value = 3
$ws.Cells.Item(98, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(98, 3).Value = "Coquimbo"
$ws.Cells.Item(98, 4).Value = 44476
$ws.Cells.Item(98, 5).Value = 5
$ws.Cells.Item(98, 6).Value = "Fruta"
$ws.Cells.Item(98, 7).Value = 100101
$ws.Cells.Item(98, 8).Value = "Berries"
$ws.Cells.Item(98, 9).Value = 100101001
$ws.Cells.Item(98, 10).Value = "Arándano (blue)"
$ws.Cells.Item(98, 11).Value = "Sin especificar"
$ws.Cells.Item(98, 12).Value = "Primera"
$ws.Cells.Item(98, 13).Value = 45
$ws.Cells.Item(98, 14).Value = 10000
$ws.Cells.Item(98, 15).Value = 10000
$ws.Cells.Item(98, 16).Value = 10000
$ws.Cells.Item(98, 17).Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(98, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(98, 19).Value = 6667
$ws.Cells.Item(98, 20).Value = 1.5
